$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(46).EntireRow.Insert()

$ws.Cells.Item(46,1).Value = 9
$ws.Cells.Item(46,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(46,3).Value = "Metropolitana"
$ws.Cells.Item(46,4).Value = 44474
$ws.Cells.Item(46,5).Value = 13
$ws.Cells.Item(46,6).Value = "Fruta"
$ws.Cells.Item(46,7).Value = 100102
$ws.Cells.Item(46,8).Value = "Cítricos"
$ws.Cells.Item(46,9).Value = 100102006
$ws.Cells.Item(46,10).Value = "Pomelo"
$ws.Cells.Item(46,11).Value = "Start Ruby"
$ws.Cells.Item(46,12).Value = "Primera"
$ws.Cells.Item(46,13).Value = 120
$ws.Cells.Item(46,14).Value = 9000
$ws.Cells.Item(46,15).Value = 9000
$ws.Cells.Item(46,16).Value = 9000
$ws.Cells.Item(46,17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(46,18).Value = "Provincia de Limarí"
$ws.Cells.Item(46,19).Value = 643
$ws.Cells.Item(46,20).Value = 14
